$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (shifting Late/heading/Outstanding one column right)
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet / tab, with K14 selected
$ws.Activate()
$ws.Range("K14").Select()
